$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "比赛Id",
    "赛季",
    "比赛日",
    "比赛时间",
    "比赛",
    "主队",
    "客队",
    "主队中文",
    "客队中文",
    "博彩公司Id",
    "博彩公司英文名称",
    "博彩公司中文名称",
    "初盘主胜赔付",
    "初盘平局赔付",
    "初盘客胜赔付",
    "即时终盘主胜赔付",
    "即时终盘平局赔付",
    "即时终盘客胜赔付",
    "初盘主胜概率",
    "初盘主胜概率",
    "初盘主胜概率",
    "即时终盘主胜概率",
    "即时终盘平局概率",
    "即时终盘客胜概率"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
